$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 4 new localization rows (2 "building_name" + 2 "building_description"
# entries) for the new basic ammo factories, at the bottom of the data range.
# Column A (keys) is populated first across all 4 rows, then column B
# (values), matching how the shared-string table fills when typed key-column
# first, value-column second.
$ws.Cells.Item(222,1).Value = "gui/hud/building_name/ammo_factory_explosive_liquid"
$ws.Cells.Item(223,1).Value = "gui/hud/building_name/ammo_factory_lowcaliber_highcaliber"
$ws.Cells.Item(224,1).Value = "gui/hud/building_description/ammo_factory_explosive_liquid"
$ws.Cells.Item(225,1).Value = "gui/hud/building_description/ammo_factory_lowcaliber_highcaliber"

$ws.Cells.Item(222,2).Value = "Explosives and liquid ammo factory"
$ws.Cells.Item(223,2).Value = "Low and High Caliber ammo factory"
$ws.Cells.Item(224,2).Value = "Produces basic explosives and liquid ammunitions"
$ws.Cells.Item(225,2).Value = "Produces basic low and ligh caliber ammunitions"

# Re-sort the whole data range (A2:K225) ascending by column A, the same way
# the existing autoFilter/sortState on the sheet already describes.
$rng = $ws.Range("A2:K225")
$rng.Sort($ws.Range("A2"))

# Restore the view: scrolled back to the top, with B14 selected.
$ws.Application.ActiveWindow.ScrollRow = 1
$sel = $ws.Range("B14")
$sel.Select()
